$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# 1) Remove the "Meta description" paragraph that currently sits right
#    after the title heading ("Play Black Hawk Deluxe Free: Review &
#    Features 2021" / "Meta description: Discover Black Hawk Deluxe...").
# -----------------------------------------------------------------------
$metaFind = $d.Content
$metaFind.Find.Execute("Meta description", $true, $false, $false, $false, `
    $false, $true, 1, $false, "", 0) | Out-Null
$metaPara = $metaFind.Paragraphs(1)
$metaPara.Range.Delete()

# -----------------------------------------------------------------------
# 2) Locate the trailing italic "image prompt" paragraph (the last
#    paragraph in the document) and insert a new bold paragraph
#    containing the title text right before it.
# -----------------------------------------------------------------------
$imgFind = $d.Content
$imgFind.Find.Execute("Create a cartoon-style feature image", $true, $false, `
    $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$imgPara = $imgFind.Paragraphs(1)
$imgParaIndex = $imgPara.Range.Paragraphs(1).Index

$insertPoint = $d.Range($imgPara.Range.Start, $imgPara.Range.Start)
$newParaXml = '<?xml version="1.0" standalone="yes"?>' + `
    '<?mso-application progid="Word.Document"?>' + `
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" ' + `
    'pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData>' + `
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body>' + `
    '<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr>' + `
    '<w:t>Play Black Hawk Deluxe Free: Review &amp; Features 2021</w:t></w:r></w:p>' + `
    '<w:p/>' + `
    '</w:body></w:document>' + `
    '</pkg:xmlData></pkg:part></pkg:package>'
$insertPoint.InsertXML($newParaXml)

# InsertXML leaves a stray empty paragraph behind (the artifact used to
# force the paragraph break) immediately after the new bold paragraph -
# remove it.
$artifactPara = $d.Paragraphs($imgParaIndex + 1)
$artifactPara.Range.Delete()

# -----------------------------------------------------------------------
# 3) Replace the text of the (still italic) image-prompt paragraph with
#    the meta-description copy, re-locating it fresh since indices/ranges
#    shifted after the insert/delete above.
# -----------------------------------------------------------------------
$finalFind = $d.Content
$finalFind.Find.Execute("Create a cartoon-style feature image", $true, `
    $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$finalPara = $finalFind.Paragraphs(1)
$finalRange = $d.Range($finalPara.Range.Start, $finalPara.Range.End)
$finalRange.Text = "Discover Black Hawk Deluxe and its unique features. " + `
    "Play this slot game for free and experience stunning graphics and " + `
    "customizable interface."
